$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Remove the stray "_GoBack" bookmark that sits right after "MP73010" in
#    the title line (it gets relocated below, closer to where Peixiao Wang's
#    most recent edit ended up).
# ---------------------------------------------------------------------------
$goBack = $d.Bookmarks("_GoBack")
$goBack.Delete()

# ---------------------------------------------------------------------------
# 2) Find the "Ben changing things up!" paragraph and add a brand-new
#    paragraph right after it: "Peixiao Wang things up!" in 22pt (sz=44)
#    text, with the proof-reader having flagged "Peixiao" as an unrecognised
#    word (spellStart/spellEnd) and the cursor's last position ("_GoBack")
#    landing right after "Wang".
# ---------------------------------------------------------------------------
$benPara = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*Ben changing things up!*") {
        $benPara = $p
    }
}

$benPara.Range.InsertParagraphAfter()
$newPara = $benPara.Next(1)

$newParaXml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
  '<w:pPr><w:rPr><w:sz w:val="44"/><w:szCs w:val="44"/></w:rPr></w:pPr>' +
  '<w:proofErr w:type="spellStart"/>' +
  '<w:r><w:rPr><w:sz w:val="44"/><w:szCs w:val="44"/></w:rPr><w:t>Peixiao</w:t></w:r>' +
  '<w:proofErr w:type="spellEnd"/>' +
  '<w:r><w:rPr><w:sz w:val="44"/><w:szCs w:val="44"/></w:rPr><w:t xml:space="preserve"> Wang</w:t></w:r>' +
  '<w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/>' +
  '<w:r><w:rPr><w:sz w:val="44"/><w:szCs w:val="44"/></w:rPr><w:t xml:space="preserve"> things up!</w:t></w:r>' +
  '</w:p>'

[void]$newPara.Range.InsertXML($newParaXml)
